# Update the "取得日時" (acquisition timestamp) column on the first sheet
# (ランサーズ) so that every data row (2-11) reflects the new run time.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-11 13:19:32"

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
